# This script updates the two-digit division answer table in place.
# Each populated table cell is addressed by its (row, column) position
# and its Range.Text is reassigned, which preserves the existing run
# formatting (font/size) defined on the paragraph's run properties.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "11÷8=1, 3"  # was: 14÷7=2, 0
$t.Cell(1,2).Range.Text = "29÷4=7, 1"  # was: 21÷7=3, 0
$t.Cell(1,3).Range.Text = "16÷3=5, 1"  # was: 61÷9=6, 7
$t.Cell(1,4).Range.Text = "32÷4=8, 0"  # was: 66÷4=16, 2
$t.Cell(1,5).Range.Text = "84÷6=14, 0"  # was: 21÷6=3, 3
$t.Cell(5,1).Range.Text = "58÷8=7, 2"  # was: 10÷2=5, 0
$t.Cell(5,2).Range.Text = "22÷8=2, 6"  # was: 83÷8=10, 3
$t.Cell(5,3).Range.Text = "74÷4=18, 2"  # was: 63÷5=12, 3
$t.Cell(5,4).Range.Text = "15÷6=2, 3"  # was: 92÷8=11, 4
$t.Cell(5,5).Range.Text = "28÷4=7, 0"  # was: 43÷7=6, 1
$t.Cell(9,1).Range.Text = "33÷6=5, 3"  # was: 30÷5=6, 0
$t.Cell(9,2).Range.Text = "13÷7=1, 6"  # was: 73÷7=10, 3
$t.Cell(9,3).Range.Text = "93÷2=46, 1"  # was: 76÷5=15, 1
$t.Cell(9,4).Range.Text = "98÷3=32, 2"  # was: 81÷4=20, 1
$t.Cell(9,5).Range.Text = "96÷5=19, 1"  # was: 29÷4=7, 1
$t.Cell(13,1).Range.Text = "84÷8=10, 4"  # was: 30÷3=10, 0
$t.Cell(13,2).Range.Text = "72÷6=12, 0"  # was: 94÷7=13, 3
$t.Cell(13,3).Range.Text = "56÷5=11, 1"  # was: 60÷3=20, 0
$t.Cell(13,4).Range.Text = "61÷7=8, 5"  # was: 80÷5=16, 0
$t.Cell(13,5).Range.Text = "33÷2=16, 1"  # was: 33÷5=6, 3
$t.Cell(17,1).Range.Text = "79÷9=8, 7"  # was: 10÷2=5, 0
$t.Cell(17,2).Range.Text = "14÷7=2, 0"  # was: 19÷9=2, 1
$t.Cell(17,3).Range.Text = "84÷5=16, 4"  # was: 24÷9=2, 6
$t.Cell(17,4).Range.Text = "74÷8=9, 2"  # was: 94÷4=23, 2
$t.Cell(17,5).Range.Text = "16÷2=8, 0"  # was: 57÷5=11, 2
